$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

$ws.Range("D3").Value = "Y"
$ws.Range("E3").Value = "SKIP"

$ws.Range("E6").Select()
